# S0154_EditCase.xlsx - refresh the three caseID test-data rows and record
# the user's last selection (A2:A4) on Sheet1, matching the "Java Pgm
# Completed and Properties file included" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new caseID values as text formulas first (a plain
# Range.Value assignment of a zero-padded numeric-looking string gets
# auto-coerced to a number, which would strip the leading zeros).
$ws.Range("A2").Formula = "=""00001324"""
$ws.Range("A3").Formula = "=""00001325"""
$ws.Range("A4").Formula = "=""00001327"""

# Collapse the formulas down to plain text values (Copy + PasteSpecial
# values-only) so the cells keep their original "General" number format/
# style instead of picking up an explicit Text-format style.
$ws.Range("A2:A4").Copy()
$ws.Range("A2:A4").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Leave the same selection the author ended up with.
$ws.Range("A2:A4").Select()
